$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking strings (e.g. "0.9988") stay as text,
# matching the original inlineStr cell type; restore style afterward.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.506.52"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.648.32"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "0.9987"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "300.43"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "0.3804"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").Value = "50.60"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "0.3514"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "1.225"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "0.08081"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "0.9989"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "22.12"
$ws.Range("D14").Value = "6.327"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "7.271"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "0.00001215"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "1.650.26"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "95.02"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").Value = "0.06973"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "6.636"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "0.9990"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "12.49"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").Value = "23.507.24"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "2.425"
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("D26").Value = "2.980"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "21.07"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "151.40"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "5.183"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "132.04"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "1.837.61"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "6.870"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").Value = "2.138"
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("D34").Value = "11.23"
$ws.Range("E34").Value = "  -8.05%  "
$ws.Range("D35").Value = "0.9922"
$ws.Range("E35").Value = "  -5.97%  "
$ws.Range("D36").Value = "0.02702"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("D37").Value = "0.08788"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "5.935"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").Value = "0.2428"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "0.06803"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").Value = "12.86"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").Value = "0.6877"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("D45").Value = "0.9985"
$ws.Range("D46").Value = "0.6393"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").Value = "2.254"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "3.923"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "0.07698"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "127.19"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "1.233"
$ws.Range("E51").Value = "  +2.62%  "

# Restore original (default) style now that values are committed as text
$ws.Range("D2:D51").Style = "Normal"
